# Refactor: update scraping workflow results and selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The scraped "Pages" value for the 2nd book row (row 3) changed from 14 to 50.
$ws.Range("E3").Value = 50

# Active selection moved from C5 to E2 before the file was saved.
$ws.Range("E2").Select()
